# Generate Report for Handoff
#
# The 81d9b0b3-ef76-46a3-98eb-92e1d16c015c file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff" status, with an
# updated handoff timestamp, and the zh-cn / de-de per-language sheets now
# carry an Error Detail message about the stale handback file.

$wb = $excel.ActiveWorkbook

$statusOld = "Handed back: in sync with en-US"
$statusNew = "Ready for handoff"
$overviewHandoffDate = "2016-08-31 02:55:21"
$zhCnHandoffDate = "2016-08-31 02:55:17"
$deDeHandoffDate = "2016-08-31 02:55:21"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/152ad5ac3904cc360c9be5c55693a6234396c9ac/e2e/81d9b0b3-ef76-46a3-98eb-92e1d16c015c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f1a7a99049a536a5cf71902c1dcd89ab872b90b/e2e/81d9b0b3-ef76-46a3-98eb-92e1d16c015c.md."

# ---- Overview sheet: row 3 is the 81d9b0b3...md file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew
$wsOverview.Range("G3").Value = $overviewHandoffDate

# ---- zh-cn sheet: row 3 is the 81d9b0b3...md file ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusNew
$wsZhCn.Range("H3").Value = $zhCnHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: row 3 is the 81d9b0b3...md file ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusNew
$wsDeDe.Range("H3").Value = $deDeHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
